$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (61) down to the new row (62)
$ws.Range("A61:V61").Copy()
$ws.Range("A62:V62").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$row = 62

$ws.Cells.Item($row, 1).Value = 61
$ws.Cells.Item($row, 2).Value = "azerbaijan"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45241.5625
$ws.Cells.Item($row, 6).Value = "Neftci Baku"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Sumqayit"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 1.5
$ws.Cells.Item($row, 11).Value = "10/11/2023 01:42"
$ws.Cells.Item($row, 12).Value = 1.54
$ws.Cells.Item($row, 13).Value = "11/11/2023 13:07"
$ws.Cells.Item($row, 14).Value = 3.76
$ws.Cells.Item($row, 15).Value = "10/11/2023 01:42"
$ws.Cells.Item($row, 16).Value = 3.84
$ws.Cells.Item($row, 17).Value = "11/11/2023 13:18"
$ws.Cells.Item($row, 18).Value = 5.5
$ws.Cells.Item($row, 19).Value = "10/11/2023 01:42"
$ws.Cells.Item($row, 20).Value = 6.17
$ws.Cells.Item($row, 21).Value = "11/11/2023 13:07"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/neftci-baku-sumqayit-fk/t8iUP5kT/"

Write-Output "done"
